$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Bulan" value (B6) used to show the month placeholder; switch it to a
# start/end date range placeholder instead.
$ws.Range("B6").Value = ": [d.start_date] – [d.end_date]"

# Reflect the edit by leaving the selection on the cell that was just
# changed (was A8, now B6).
$ws.Range("B6").Select()
